$wb = $excel.ActiveWorkbook

# --- Sheet: ALC (index 1) ---
$ws = $wb.Worksheets.Item(1)
$ws.Range("J105").Value = 81000
$ws.Range("L105").Value = 81000
$ws.Range("H105").Value = 81000
$ws.Range("N105").Value = -87988
$ws.Range("J112").Value = 5000
$ws.Range("N112").Value = -17216
$ws.Range("L112").Value = 15000
$ws.Range("H112").Value = 4350

# --- Sheet: ARM (index 2) ---
$ws = $wb.Worksheets.Item(2)
$ws.Range("I2").Value = 499
$ws.Range("M2").Value = -386
$ws.Range("K2").Value = 499
$ws.Range("H2").Value = 646
$ws.Range("J2").Value = 851.8
$ws.Range("N2").Value = -1077.8
$ws.Range("L2").Value = 851.8
$ws.Range("H45").Value = 2803.65
$ws.Range("K45").Value = 1845.5
$ws.Range("I45").Value = 1845.5
$ws.Range("M45").Value = -1468.5
$ws.Range("H74").Value = 787.65625
$ws.Range("N74").Value = -2602.2
$ws.Range("M74").Value = 197.25
$ws.Range("L74").Value = 854.2
$ws.Range("K74").Value = 676.75
$ws.Range("J74").Value = 854.2
$ws.Range("I74").Value = 676.75
$ws.Range("H77").Value = 787.65625
$ws.Range("N77").Value = -13007
$ws.Range("K77").Value = 3383.75
$ws.Range("L77").Value = 4271
$ws.Range("I77").Value = 676.75
$ws.Range("M77").Value = 984.25
$ws.Range("J77").Value = 854.2
$ws.Range("I97").Value = 787.65717
$ws.Range("N97").Value = -2245.1428
$ws.Range("H97").Value = 920.6531
$ws.Range("L97").Value = 1253.1428
$ws.Range("M97").Value = -291.65717
$ws.Range("J97").Value = 1253.1428
$ws.Range("K97").Value = 787.65717
$ws.Range("J101").Value = 100602
$ws.Range("N101").Value = -107092
$ws.Range("H101").Value = 100602
$ws.Range("L101").Value = 100602
$ws.Range("L102").Value = 3000
$ws.Range("K102").Value = 1562
$ws.Range("M102").Value = 60
$ws.Range("J102").Value = 3000
$ws.Range("H102").Value = 1801.6666
$ws.Range("I102").Value = 1562
$ws.Range("H113").Value = 0
$ws.Range("J113").Value = 0
$ws.Range("L113").Value = 0
$ws.Range("L114").Value = 0
$ws.Range("H114").Value = 0
$ws.Range("J114").Value = 0
$ws.Range("M116").Value = 1795
$ws.Range("K116").Value = 499
$ws.Range("J116").Value = 851.8
$ws.Range("N116").Value = -5439.8
$ws.Range("I116").Value = 499
$ws.Range("H116").Value = 646
$ws.Range("L116").Value = 851.8
$ws.Range("J117").Value = 100248
$ws.Range("N117").Value = -109426
$ws.Range("L117").Value = 100248
$ws.Range("H117").Value = 100248
$ws.Range("L119").Value = 0
$ws.Range("J119").Value = 0
$ws.Range("H119").Value = 0
$ws.Range("N102").Value = -6244
$ws.Range("N113").ClearContents()
$ws.Range("N114").ClearContents()
$ws.Range("N119").ClearContents()

# --- Sheet: BSM (index 3) ---
$ws = $wb.Worksheets.Item(3)
$ws.Range("I3").Value = 499
$ws.Range("M3").Value = -385
$ws.Range("L3").Value = 851.8
$ws.Range("H3").Value = 646
$ws.Range("N3").Value = -1079.8
$ws.Range("J3").Value = 851.8
$ws.Range("K3").Value = 499
$ws.Range("J92").Value = 0
$ws.Range("L92").Value = 0
$ws.Range("H92").Value = 0
$ws.Range("J104").Value = 0
$ws.Range("H104").Value = 0
$ws.Range("L104").Value = 0
$ws.Range("J105").Value = 3470.3333
$ws.Range("M105").Value = -316.3332999999998
$ws.Range("K105").Value = 2063.3333
$ws.Range("L105").Value = 3470.3333
$ws.Range("I105").Value = 2063.3333
$ws.Range("H105").Value = 2532.3333
$ws.Range("N105").Value = -6964.3333
$ws.Range("H110").Value = 0
$ws.Range("J110").Value = 0
$ws.Range("L110").Value = 0
$ws.Range("L111").Value = 0
$ws.Range("J111").Value = 0
$ws.Range("H111").Value = 0
$ws.Range("J117").Value = 0
$ws.Range("L117").Value = 0
$ws.Range("H117").Value = 0
$ws.Range("N92").ClearContents()
$ws.Range("N104").ClearContents()
$ws.Range("N110").ClearContents()
$ws.Range("N111").ClearContents()
$ws.Range("N117").ClearContents()

# --- Sheet: CRP (index 4) ---
$ws = $wb.Worksheets.Item(4)
$ws.Range("N41").Value = -18982.666
$ws.Range("L41").Value = 18126.666
$ws.Range("J41").Value = 18126.666
$ws.Range("H41").Value = 18126.666
$ws.Range("J43").Value = 95000
$ws.Range("N43").Value = -95368
$ws.Range("L43").Value = 95000
$ws.Range("H43").Value = 95000
$ws.Range("I58").Value = 2088.6843
$ws.Range("J58").Value = 876
$ws.Range("N58").Value = -1282
$ws.Range("K58").Value = 2088.6843
$ws.Range("H58").Value = 1729.3704
$ws.Range("L58").Value = 876
$ws.Range("M58").Value = -1885.6843
$ws.Range("J95").Value = 100624
$ws.Range("H95").Value = 100624
$ws.Range("L95").Value = 100624
$ws.Range("H96").Value = 89812
$ws.Range("N96").Value = -95304
$ws.Range("J96").Value = 89812
$ws.Range("L96").Value = 89812
$ws.Range("J101").Value = 95000
$ws.Range("N101").Value = -101490
$ws.Range("H101").Value = 95000
$ws.Range("L101").Value = 95000
$ws.Range("M105").Value = 1416.8889
$ws.Range("K105").Value = 330.1111
$ws.Range("I105").Value = 330.1111
$ws.Range("H105").Value = 330.1111
$ws.Range("J107").Value = 1675
$ws.Range("K107").Value = 4166975
$ws.Range("M107").Value = -4165055
$ws.Range("L107").Value = 1675
$ws.Range("H107").Value = 3290069.8
$ws.Range("N107").Value = -5515
$ws.Range("I107").Value = 4166975
$ws.Range("N111").Value = -108882
$ws.Range("L111").Value = 100702
$ws.Range("J111").Value = 100702
$ws.Range("H111").Value = 100702
$ws.Range("L114").Value = 98684
$ws.Range("H114").Value = 98684
$ws.Range("N114").Value = -107362
$ws.Range("J114").Value = 98684
$ws.Range("H115").Value = 64645
$ws.Range("J115").Value = 64645
$ws.Range("L115").Value = 64645
$ws.Range("J116").Value = 98580.664
$ws.Range("N116").Value = -107758.664
$ws.Range("H116").Value = 98580.664
$ws.Range("L116").Value = 98580.664
$ws.Range("J117").Value = 99712
$ws.Range("N117").Value = -108890
$ws.Range("L117").Value = 99712
$ws.Range("H117").Value = 99712
$ws.Range("L118").Value = 82075.336
$ws.Range("H118").Value = 82075.336
$ws.Range("J118").Value = 82075.336
$ws.Range("L119").Value = 0
$ws.Range("J119").Value = 0
$ws.Range("H119").Value = 0
$ws.Range("K136").Value = 6266.0529
$ws.Range("N136").Value = -7728
$ws.Range("H136").Value = 1729.3704
$ws.Range("I136").Value = 2088.6843
$ws.Range("L136").Value = 2628
$ws.Range("J136").Value = 876
$ws.Range("M136").Value = -3716.0529
$ws.Range("N95").Value = -106116
$ws.Range("N115").Value = -66995
$ws.Range("N118").Value = -85389.336
$ws.Range("N119").ClearContents()

# --- Sheet: CUL (index 5) ---
$ws = $wb.Worksheets.Item(5)
$ws.Range("H48").Value = 12990
$ws.Range("L48").Value = 38970
$ws.Range("N48").Value = -39470
$ws.Range("J48").Value = 12990
$ws.Range("H68").Value = 790.71
$ws.Range("L68").Value = 2533.8261
$ws.Range("J68").Value = 844.6087
$ws.Range("K68").Value = 2012.22582
$ws.Range("M68").Value = -1201.22582
$ws.Range("N68").Value = -4155.8261
$ws.Range("I68").Value = 670.74194
$ws.Range("L71").Value = 7601.4783
$ws.Range("J71").Value = 844.6087
$ws.Range("N71").Value = -15713.4783
$ws.Range("M71").Value = -1980.67746
$ws.Range("H71").Value = 790.71
$ws.Range("I71").Value = 670.74194
$ws.Range("K71").Value = 6036.67746
$ws.Range("H113").Value = 789.9286
$ws.Range("I113").Value = 454
$ws.Range("J113").Value = 1081.0667
$ws.Range("M113").Value = 808
$ws.Range("K113").Value = 1362
$ws.Range("L113").Value = 3243.2001
$ws.Range("N113").Value = -7583.2001
$ws.Range("M122").Value = -430
$ws.Range("H122").Value = 5302.048
$ws.Range("I122").Value = 320
$ws.Range("K122").Value = 2880
$ws.Range("J134").Value = 6666.6665
$ws.Range("M134").Value = -300
$ws.Range("H134").Value = 2704.375
$ws.Range("I134").Value = 1790
$ws.Range("L134").Value = 19999.9995
$ws.Range("K134").Value = 5370
$ws.Range("N134").Value = -30139.9995
$ws.Range("I139").Value = 1283.3334
$ws.Range("K139").Value = 3850.0002
$ws.Range("M139").Value = 1289.9998
$ws.Range("N139").Value = -26280.0005
$ws.Range("J139").Value = 5333.3335
$ws.Range("H139").Value = 3308.3333
$ws.Range("L139").Value = 16000.0005

# --- Sheet: GSM (index 6) ---
$ws = $wb.Worksheets.Item(6)
$ws.Range("L80").Value = 143002
$ws.Range("K80").Value = 4504000
$ws.Range("M80").Value = -4503002
$ws.Range("H80").Value = 1887401.2
$ws.Range("J80").Value = 143002
$ws.Range("N80").Value = -144998
$ws.Range("I80").Value = 4504000
$ws.Range("H83").Value = 1887401.2
$ws.Range("N83").Value = -724994
$ws.Range("L83").Value = 715010
$ws.Range("I83").Value = 4504000
$ws.Range("M83").Value = -22515008
$ws.Range("J83").Value = 143002
$ws.Range("K83").Value = 22520000
$ws.Range("J101").Value = 0
$ws.Range("H101").Value = 0
$ws.Range("L101").Value = 0
$ws.Range("M122").Value = -1536.5716
$ws.Range("H122").Value = 2144.6667
$ws.Range("I122").Value = 1328.8572
$ws.Range("K122").Value = 3986.5716
$ws.Range("I132").Value = 2332.4614
$ws.Range("N132").Value = -14157.2
$ws.Range("L132").Value = 9097.200000000001
$ws.Range("H132").Value = 2588.5366
$ws.Range("K132").Value = 6997.3842
$ws.Range("M132").Value = -4467.3842
$ws.Range("J132").Value = 3032.4
$ws.Range("N101").ClearContents()

# --- Sheet: LTW (index 7) ---
$ws = $wb.Worksheets.Item(7)
$ws.Range("K93").Value = 7727.9287
$ws.Range("L93").Value = 3262.4
$ws.Range("N93").Value = -5758.4
$ws.Range("J93").Value = 3262.4
$ws.Range("I93").Value = 7727.9287
$ws.Range("M93").Value = -6479.9287
$ws.Range("H93").Value = 6552.7896

# --- Sheet: WVR (index 8) ---
$ws = $wb.Worksheets.Item(8)
$ws.Range("M122").Value = -625.4998000000001
$ws.Range("H122").Value = 1539.5555
$ws.Range("I122").Value = 1025.1666
$ws.Range("L122").Value = 7704.999899999999
$ws.Range("K122").Value = 3075.4998
$ws.Range("J122").Value = 2568.3333
$ws.Range("N122").Value = -12604.9999
